$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.362.90"
$ws.Range("E2").Value = "  -3.44%  "

# Row 3
$ws.Range("D3").Value = "1.648.93"
$ws.Range("E3").Value = "  -3.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "

# Row 6
$ws.Range("E6").Value = "  -2.70%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0615"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.50%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("D12").Value = "1.884.54"
$ws.Range("E12").Value = "  -3.57%  "

# Row 13
$ws.Range("D13").Value = "1.654.01"
$ws.Range("E13").Value = "  -3.37%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.34%  "

# Row 17
$ws.Range("D17").Value = "27.358.95"
$ws.Range("E17").Value = "  -3.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.53%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  -2.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.41%  "

# Row 21
$ws.Range("E21").Value = "  +0.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.99%  "

# Row 24
$ws.Range("E24").Value = "  -0.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "

# Row 26
$ws.Range("E26").Value = "  -2.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.37%  "

# Row 28
$ws.Range("E28").Value = "  +0.09%  "

# Row 29
$ws.Range("E29").Value = "  -2.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.85%  "

# Row 31
$ws.Range("E31").Value = "  -1.19%  "

# Row 32
$ws.Range("E32").Value = "  -2.83%  "

# Row 33
$ws.Range("D33").Value = "1.458.82"
$ws.Range("E33").Value = "  -1.54%  "

# Row 34
$ws.Range("E34").Value = "  -3.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.31%  "

# Row 36
$ws.Range("E36").Value = "  -0.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.909"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.59%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.67%  "

# Row 40
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.23%  "

# Row 44
$ws.Range("E44").Value = "  -1.96%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.791.99"
$ws.Range("E45").Value = "  -3.63%  "

# Row 46
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.63%  "

# Row 47
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.13%  "

# Row 49
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("E50").Value = "  -2.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.14%  "
